# Apply updated cryptocurrency price/volume data scraped on 2023-09-30.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.030.55'
$ws.Range("E2").Value = '  +0.78%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.679.55'
$ws.Range("E3").Value = '  +1.08%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.87'
$ws.Range("E5").Value = '  +0.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.519'
$ws.Range("E6").Value = '  -2.89%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.254'
$ws.Range("E8").Value = '  +1.70%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '21.38'
$ws.Range("E9").Value = '  +6.08%  '

$ws.Range("E10").Value = '  +0.63%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0888'
$ws.Range("E11").Value = '  -0.81%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.917.85'
$ws.Range("E12").Value = '  +1.18%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.670.87'
$ws.Range("E13").Value = '  +0.56%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.10'
$ws.Range("E14").Value = '  +0.56%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.532'
$ws.Range("E15").Value = '  +1.83%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.46'
$ws.Range("E16").Value = '  +0.74%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.040.98'
$ws.Range("E17").Value = '  +0.83%  '

$ws.Range("E18").Value = '  +4.30%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '235.02'
$ws.Range("E19").Value = '  +1.30%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0737'
$ws.Range("E20").Value = '  +0.97%  '

$ws.Range("E21").Value = '  -0.01%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.46'
$ws.Range("E22").Value = '  +0.90%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.27'
$ws.Range("E23").Value = '  +1.23%  '

$ws.Range("E24").Value = '  -4.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.40'
$ws.Range("E25").Value = '  +0.51%  '

$ws.Range("E26").Value = '  +1.75%  '

$ws.Range("E27").Value = '  +3.62%  '

$ws.Range("E28").Value = '  -2.34%  '

$ws.Range("E29").Value = '  +0.05%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0498'
$ws.Range("E30").Value = '  +0.53%  '

$ws.Range("E31").Value = '  +0.28%  '

$ws.Range("E32").Value = '  +0.60%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.536.76'
$ws.Range("E33").Value = '  +5.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.17'
$ws.Range("E34").Value = '  +0.72%  '

$ws.Range("E35").Value = '  +4.99%  '

$ws.Range("E36").Value = '  -0.65%  '

$ws.Range("E37").Value = '  +3.00%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.922'
$ws.Range("E38").Value = '  +2.81%  '

$ws.Range("E39").Value = '  +3.41%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.04'
$ws.Range("E40").Value = '  +6.28%  '

$ws.Range("E41").Value = '  -0.01%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.69'
$ws.Range("E42").Value = '  -2.72%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '67.90'
$ws.Range("E43").Value = '  +3.33%  '

$ws.Range("E44").Value = '  -0.12%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.822.13'
$ws.Range("E45").Value = '  +0.58%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.781'
$ws.Range("E46").Value = '  +0.53%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.07'
$ws.Range("E47").Value = '  -0.35%  '

$ws.Range("E48").Value = '  +0.12%  '

$ws.Range("E49").Value = '  +2.57%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.00'
$ws.Range("E50").Value = '  +5.80%  '
